$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

# Mark column B as "OK" for the rows whose review text has now been
# fully incorporated (V1 da Dissertacao com o texto completo).
$rows = @(2, 8, 16, 19, 21, 22, 25, 26)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 2).Value = "OK"
}

# Update the view: move the active selection to B27.
$ws.Activate()
$ws.Range("B27").Select()
